$wb = $excel.ActiveWorkbook

# The "想去人数" (wanted-to-go count) column F needs updating in both the
# "展览" sheet and the "全部类型" sheet, which carry the same rows.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 901
    $ws.Range("F3").Value = 4587
    $ws.Range("F4").Value = 135
    $ws.Range("F5").Value = 798
}
